# "2.xlsx" (Threshold/Zn/3His) was re-uploaded with a couple of the
# threshold numbers tweaked on Sheet1, and the selection cursor ended up
# parked on C4 when it was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Max column (C) threshold updates for alpha_distance_range / beta_distance_range / ratio_threshold_range
$ws.Range("C2").Value = 10.7               # was 11
$ws.Range("C3").Value = 9.3000000000000007 # was 9.5
$ws.Range("C4").Value = 1.45               # was 1.4

# The workbook was last saved with the cursor on C4 (previously F4).
$ws.Range("C4").Select()
